# Weekly refresh of the "Hortaliza, Vega Monumental Concepción - Betarraga"
# data block (rows 180-205): a new week (row pair 180/181) is inserted at the
# top of this market's table, every existing week shifts down by one row
# pair, and the oldest week (the old rows 204/205) is appended as new rows
# 206/207.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Append the two new rows at the bottom (206/207), copying the current
#    (soon to be shifted-away) oldest week, rows 204/205, verbatim.
$ws.Range("A204:R204").Copy($ws.Range("A206:R206"))
$ws.Range("A205:R205").Copy($ws.Range("A207:R207"))

# 2) Shift rows 182..205 down by one row pair: new row r gets the old
#    contents of row r-2. Walk from the bottom up so a source row is never
#    clobbered before it has been read.
for ($r = 205; $r -ge 182; $r--) {
    $src = $r - 2
    $ws.Range("A" + $src + ":R" + $src).Copy($ws.Range("A" + $r + ":R" + $r))
}

# 3) The freed-up top rows (180/181) become the new week: same
#    Primera/Segunda price data as before, only the date (column D) moves
#    forward to the new reporting date.
$ws.Range("D180").Value = 44474
$ws.Range("D181").Value = 44474
